# Update "想去人数" (interested-people count) figures that changed between
# data pulls, on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F10").Value = 6890
$wsExhibit.Range("F14").Value = 7797
$wsExhibit.Range("F17").Value = 5028
$wsExhibit.Range("F19").Value = 2281
$wsExhibit.Range("F26").Value = 282
$wsExhibit.Range("F37").Value = 2094

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 6890
$wsAll.Range("F18").Value = 7797
$wsAll.Range("F21").Value = 5028
$wsAll.Range("F23").Value = 2281
$wsAll.Range("F32").Value = 282
$wsAll.Range("F44").Value = 2094
